$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(6, 1).Value = 42607.889189814814
$ws.Cells.Item(6, 2).Value = -8
$ws.Cells.Item(6, 3).Value = 57
$ws.Cells.Item(6, 4).Value = 39
$ws.Cells.Item(6, 5).Value = 31
$ws.Cells.Item(6, 6).Value = 68
$ws.Cells.Item(6, 7).Value = 38807
$ws.Cells.Item(6, 8).Value = 31379
$ws.Cells.Item(6, 9).Value = 3781
$ws.Cells.Item(6, 10).Value = 380
$ws.Cells.Item(6, 11).Value = 261
$ws.Cells.Item(6, 12).Value = 14
$ws.Cells.Item(6, 13).Value = 31
$ws.Cells.Item(6, 14).Value = "Bag"
